# CANMFA_Anschluss.xlsx edit:
# "make small text display more self-adjusting due to measured adc values"
#
# Tabelle2 (sheet2) gets a new small header/value table in column A/C
# starting at row 32: A32 = "MFA", C32 = "BB", followed by 26 measured
# ADC values in C33:C58 counting down from 26 to 1 (one per existing
# row of the KL30/KL15/... list that already occupies A33:B58).
# The user's selection ends up on the newly entered C33 cell, and
# Tabelle2 becomes the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Tabelle2")

# New small header row for the measured-ADC-value column.
$ws2.Range("A32").Value = "MFA"
$ws2.Range("C32").Value = "BB"

# Fill in the measured ADC values, counting down from 26 to 1.
For ($i = 0; $i -lt 26; $i++) {
    $row = 33 + $i
    $value = 26 - $i
    $ws2.Cells.Item($row, 3).Value = $value
}

# Reflect the author's final selection/active sheet state: Tabelle2
# becomes the active tab, with the newly entered C33 cell selected.
[void]$ws2.Activate()
[void]$ws2.Range("C33").Select()
